# feat: add 2022-Q1 data
#
# The workbook has 4 sheets: 2021-Q1, 2021-Q2, 2021-Q3, 总计.
# We need to:
#   1. Turn the old "总计" sheet into the new "2022-Q1" sheet (same
#      sheetId/rId slot) holding the Q1-2022 fund holdings detail.
#   2. Add a brand new "总计" sheet at the end with the summary table
#      updated to include the 2022-Q1 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Step 1: repurpose the existing "总计" sheet as "2022-Q1"
# ---------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Drop to just the header row + one data row (existing sheet has 4 rows).
$q1.Rows.Item(4).Delete()
$q1.Rows.Item(3).Delete()

# Header row (B1:D1 already carry the bordered/bold "s=2" style from the
# old 总计 header; extend that same style across E1:H1).
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"

$q1.Range("B1:D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data row 2 - fund holding detail for 2022-Q1.
$q1.Range("A2").Value = 0

$q1.Range("B2").Value = "'090011"
$q1.Range("B2").Style = "Normal"

$q1.Range("C2").Value = "大成核心双动力混合"
$q1.Range("C2").Style = "Normal"

$q1.Range("D2").Value = "'0.34"
$q1.Range("D2").Style = "Normal"

$q1.Range("E2").Value = "'93.14"
$q1.Range("E2").Style = "Normal"

$q1.Range("F2").Value = "'2.02"
$q1.Range("F2").Style = "Normal"

$q1.Range("G2").Value = "'0.0069"
$q1.Range("G2").Style = "Normal"

$q1.Range("H2").Value = 9

# ---------------------------------------------------------------
# Step 2: add the new "总计" sheet at the end of the tab strip
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$total = $wb.Worksheets.Add($null, $lastSheet)
$total.Name = "总计"

# Reuse the bordered/bold header style and the A-column index style from
# the "2022-Q1" sheet (same visual style the old "总计" sheet used).
$q1.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$q1.Range("A2").Copy()
$total.Range("A2:A5").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "'2022-Q1"
$total.Range("B2").Style = "Normal"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

$total.Range("A3").Value = 1
$total.Range("B3").Value = "'2021-Q3"
$total.Range("B3").Style = "Normal"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.03

$total.Range("A4").Value = 2
$total.Range("B4").Value = "'2021-Q2"
$total.Range("B4").Style = "Normal"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.02

$total.Range("A5").Value = 3
$total.Range("B5").Value = "'2021-Q1"
$total.Range("B5").Style = "Normal"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.58
